$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Add the missing journal entry on row 21
$ws.Range("A21").Value2 = 45490
$ws.Range("B21").Value2 = 0.79166666666666663
$ws.Range("C21").Value2 = 0.99930555555555556
$ws.Range("E21").Value = "Analyse et état de l'art"

$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("B20:C20").Copy()
$ws.Range("B21").PasteSpecial(-4122)

# Update sheet view: scroll / selection
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("H20").Select()

$wb.Application.Calculate()
